$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 111541119
$ws.Range("B9").Value = 5426
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 101410
$ws.Range("F9").Value = "Reliktbock"
$ws.Range("G9").Value = "Nothorhina muricata"
$ws.Range("H9").Value = "(Dalman, 1817)"
$ws.Range("M9").Value = "äldre gnagspår"
$ws.Range("Q9").Value = 693467.6220677271
$ws.Range("R9").Value = 6551532.561666255
$ws.Range("AC9").Value = "En gammal tall med kläckhål här och var. Om det är färskt eller gammalt är svårt sia om."
$ws.Range("AO9").Value = "gammeltall"

# Row 10
$ws.Range("A10").Value = 111541122
$ws.Range("B10").Value = 5112
$ws.Range("E10").Value = 102204
$ws.Range("F10").Value = "Grönhjon"
$ws.Range("G10").Value = "Callidium aeneum"
$ws.Range("H10").Value = "(De Geer, 1775)"
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").Value = "äldre gnagspår"
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 693344.0451535647
$ws.Range("R10").Value = 6551526.82974836
$ws.Range("AF10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111541128
$ws.Range("B11").Value = 5113
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 100526
$ws.Range("F11").Value = "Bronshjon"
$ws.Range("G11").Value = "Callidium coriaceum"
$ws.Range("H11").Value = "Paykull, 1800"
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = "färska gnagspår"
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value = 693570.8046739453
$ws.Range("R11").Value = 6551451.742365629
$ws.Range("AF11").ClearContents()
$ws.Range("AO11").Value = "torrgran"

# Row 12
$ws.Range("A12").Value = 111541115
$ws.Range("B12").Value = 89405
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 1202
$ws.Range("F12").Value = "Ullticka"
$ws.Range("G12").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H12").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 693612.9253791923
$ws.Range("R12").Value = 6551435.326171798
$ws.Range("AO12").Value = "låga av gran"

# Row 13
$ws.Range("A13").Value = 111541120
$ws.Range("B13").Value = 79444
$ws.Range("E13").Value = 1049
$ws.Range("F13").Value = "Kortskaftad ärgspik"
$ws.Range("G13").Value = "Microcalicium ahlneri"
$ws.Range("H13").Value = "Tibell"
$ws.Range("J13").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("Q13").Value = 693513.2669972532
$ws.Range("R13").Value = 6551517.868690074
$ws.Range("AC13").ClearContents()
$ws.Range("AF13").ClearContents()
$ws.Range("AO13").Value = "silverstubbe av tall"

# Row 14
$ws.Range("A14").Value = 111541118
$ws.Range("B14").Value = 94851
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 2569
$ws.Range("F14").Value = "Stor revmossa"
$ws.Range("G14").Value = "Bazzania trilobata"
$ws.Range("H14").Value = "(L.) Gray"
$ws.Range("J14").ClearContents()
$ws.Range("K14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("Q14").Value = 693461.6376634488
$ws.Range("R14").Value = 6551559.049034445
$ws.Range("AF14").ClearContents()
$ws.Range("AO14").ClearContents()

# Row 15
$ws.Range("A15").Value = 111541121
$ws.Range("B15").Value = 79444
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1049
$ws.Range("F15").Value = "Kortskaftad ärgspik"
$ws.Range("G15").Value = "Microcalicium ahlneri"
$ws.Range("H15").Value = "Tibell"
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("Q15").Value = 693460.9606228607
$ws.Range("R15").Value = 6551521.405726598
$ws.Range("AF15").ClearContents()
$ws.Range("AO15").Value = "silverstubbe av tall"

# Row 16
$ws.Range("A16").Value = 111541129
$ws.Range("B16").Value = 5113
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 100526
$ws.Range("F16").Value = "Bronshjon"
$ws.Range("G16").Value = "Callidium coriaceum"
$ws.Range("H16").Value = "Paykull, 1800"
$ws.Range("J16").ClearContents()
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = "äldre gnagspår"
$ws.Range("N16").ClearContents()
$ws.Range("Q16").Value = 693328.6441019299
$ws.Range("R16").Value = 6551545.628735202
$ws.Range("AF16").ClearContents()
$ws.Range("AO16").Value = "silverstubbe av tall"
